# The deck currently has the "Integral" theme colours applied to its
# (only) slide master / live theme part. The commit swaps the presentation
# back to the stock "Office Theme" colour scheme (the other theme part in
# the package, which is otherwise only reachable from the Notes Master,
# already carries the "Office Theme" name/colours/fonts/formats).
#
# Re-apply the 12 "Office Theme" colour-scheme entries (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) onto the presentation's active theme colour
# scheme, in MsoThemeColorSchemeIndex order, via the Slide's ThemeColorScheme
# (which maps straight onto <a:clrScheme> of the theme part backing the
# slide master that every slide/layout inherits from).

$p = $ppt.ActivePresentation

# Office Theme palette (RRGGBB) in dk1,lt1,dk2,lt2,accent1-6,hlink,folHlink order
$officeThemeRGB = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

function ToComRGB($rrggbb) {
    $r = ($rrggbb -band 0xFF0000) -shr 16
    $g = ($rrggbb -band 0x00FF00) -shr 8
    $b = ($rrggbb -band 0x0000FF)
    return $r -bor ($g -shl 8) -bor ($b -shl 16)
}

$tcs = $p.Slides.Item(1).ThemeColorScheme

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = ToComRGB($officeThemeRGB[$i - 1])
}
